$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (lamda_1) and Column C (lamda_2) are constant across all data rows (2-55).
$ws.Range("B2:B55").Value = 33.94444444444444
$ws.Range("C2:C55").Value = 1.95

# Column D (dic_nbre_clients_poisson_2_keys) and Column E (dic_nbre_clients_prob_poisson_2_values)
# new values per row, taken from the updated workbook.
$rows = @(
    @(2, 0, 0.135),
    @(3, 2, 0.001),
    @(4, 3, 0.004),
    @(5, 4, 0.006),
    @(6, 5, 0.016),
    @(7, 6, 0.032),
    @(8, 7, 0.042),
    @(9, 8, 0.043),
    @(10, 9, 0.056),
    @(11, 10, 0.036),
    @(12, 11, 0.028),
    @(13, 12, 0.028),
    @(14, 13, 0.022),
    @(15, 14, 0.028),
    @(16, 15, 0.031),
    @(17, 16, 0.037),
    @(18, 17, 0.039),
    @(19, 18, 0.029),
    @(20, 19, 0.04),
    @(21, 20, 0.029),
    @(22, 21, 0.03),
    @(23, 22, 0.022),
    @(24, 23, 0.02),
    @(25, 24, 0.02),
    @(26, 25, 0.025),
    @(27, 26, 0.019),
    @(28, 27, 0.016),
    @(29, 28, 0.02),
    @(30, 29, 0.015),
    @(31, 30, 0.017),
    @(32, 31, 0.014),
    @(33, 32, 0.018),
    @(34, 33, 0.01),
    @(35, 34, 0.012),
    @(36, 35, 0.008),
    @(37, 36, 0.009000000000000001),
    @(38, 37, 0.004),
    @(39, 38, 0.003),
    @(40, 39, 0.009000000000000001),
    @(41, 40, 0.003),
    @(42, 41, 0.002),
    @(43, 42, 0.005),
    @(44, 43, 0.004),
    @(45, 44, 0.001),
    @(46, 45, 0.001),
    @(47, 46, 0.001),
    @(48, 48, 0.001),
    @(49, 49, 0.001),
    @(50, 50, 0.001),
    @(51, 52, 0.001),
    @(52, 54, 0.002),
    @(53, 55, 0.001),
    @(54, 60, 0.001),
    @(55, 65, 0.001)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $dVal = $r[1]
    $eVal = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $dVal
    $ws.Cells.Item($rowNum, 5).Value = $eVal
}
